$p = $ppt.ActivePresentation

function Get-SlideTitleText($slide) {
    if ($slide.Shapes.HasTitle) {
        return $slide.Shapes.Title.TextFrame.TextRange.Text
    }
    return ""
}

function Find-SlideByTitle($title) {
    for ($i = 1; $i -le $p.Slides.Count; $i++) {
        $candidate = $p.Slides.Item($i)
        if ((Get-SlideTitleText $candidate) -eq $title) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Slide "Project Document": fix up the second subtitle paragraph text
#    "Updated sections" -> "Updated a few sections"
# ---------------------------------------------------------------------------
$projectDocSlide = Find-SlideByTitle "Project Document"
$subtitleShape = $projectDocSlide.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange
$secondPara = $subtitleRange.Paragraphs(2, 1)
$secondPara.Runs(1).Text = "Updated a few sections"

# ---------------------------------------------------------------------------
# 2) Delete the "Project Plan" slide entirely
# ---------------------------------------------------------------------------
$projectPlanSlide = Find-SlideByTitle "Project Plan"
$projectPlanSlide.Delete()

# ---------------------------------------------------------------------------
# 3) Duplicate the "Metronome" slide to create the new "Scaling" slide, then
#    rename the title of the original copy to "Scaling" while the duplicate
#    keeps the "Metronome" title (mirrors a "duplicate slide, then retitle"
#    authoring flow, keeping shape names/placeholders intact).
# ---------------------------------------------------------------------------
$metronomeSlide = Find-SlideByTitle "Metronome"
$metronomeIndex = $metronomeSlide.SlideIndex
$null = $metronomeSlide.Duplicate()

$scalingSlide = $p.Slides.Item($metronomeIndex)
$scalingSlide.Shapes.Title.TextFrame.TextRange.Text = "Scaling"

# ---------------------------------------------------------------------------
# 4) Move "Tuning (the beast)" to the very end, after "Scaling" and
#    "Metronome", giving the final running order:
#    ... Plan b, Scaling, Metronome, Tuning (the beast)
# ---------------------------------------------------------------------------
$tuningSlide = Find-SlideByTitle "Tuning (the beast)"
$tuningSlide.MoveTo($p.Slides.Count)
